$wb = $excel.ActiveWorkbook

# ALC row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 111
$ws.Range("I12").Value = 111
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 111
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 59
$ws.Range("N12").Value = $null

# ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1846.4117
$ws.Range("I15").Value = 1846.4117
$ws.Range("K15").Value = 5539.2351
$ws.Range("M15").Value = -5370.2351

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 283.5
$ws.Range("I33").Value = 303.8889
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 303.8889
$ws.Range("L33").Value = 100
$ws.Range("M33").Value = -74.88889999999998
$ws.Range("N33").Value = -558

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 21499.5
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").Value = $null

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 7895.0625
$ws.Range("I74").Value = 4349.9165
$ws.Range("J74").Value = 18530.5
$ws.Range("K74").Value = 4349.9165
$ws.Range("L74").Value = 18530.5
$ws.Range("M74").Value = -3413.9165
$ws.Range("N74").Value = -20402.5

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 7895.0625
$ws.Range("I77").Value = 4349.9165
$ws.Range("J77").Value = 18530.5
$ws.Range("K77").Value = 21749.5825
$ws.Range("L77").Value = 92652.5
$ws.Range("M77").Value = -17069.5825
$ws.Range("N77").Value = -102012.5

# ALC row 81
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 100000
$ws.Range("J81").Value = 100000
$ws.Range("L81").Value = 100000
$ws.Range("N81").Value = -101996

# ALC row 84
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H84").Value = 100000
$ws.Range("J84").Value = 100000
$ws.Range("L84").Value = 300000
$ws.Range("N84").Value = -309984

# ALC row 97
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 25634.334
$ws.Range("J97").Value = 32694.143
$ws.Range("L97").Value = 98082.429
$ws.Range("N97").Value = -99074.429

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 10187
$ws.Range("I113").Value = 8254.556
$ws.Range("J113").Value = 12671.571
$ws.Range("K113").Value = 8254.556
$ws.Range("L113").Value = 12671.571
$ws.Range("M113").Value = -5000.556
$ws.Range("N113").Value = -19179.571

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1673.425
$ws.Range("I132").Value = 1621.2972
$ws.Range("J132").Value = 2316.3333
$ws.Range("K132").Value = 4863.8916
$ws.Range("L132").Value = 6948.999899999999
$ws.Range("M132").Value = -2333.8916
$ws.Range("N132").Value = -12008.9999

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1608.7693
$ws.Range("I137").Value = 1563.4
$ws.Range("J137").Value = 1656.5264
$ws.Range("K137").Value = 4690.200000000001
$ws.Range("L137").Value = 4969.5792
$ws.Range("M137").Value = -2140.200000000001
$ws.Range("N137").Value = -10069.5792

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J138").Value = 8067680
$ws.Range("L138").Value = 24203040
$ws.Range("N138").Value = -24213320

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 478.72726
$ws.Range("I2").Value = 246.84616
$ws.Range("J2").Value = 813.6667
$ws.Range("K2").Value = 246.84616
$ws.Range("L2").Value = 813.6667
$ws.Range("M2").Value = -133.84616
$ws.Range("N2").Value = -1039.6667

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7493.442
$ws.Range("I32").Value = 3044.4868
$ws.Range("K32").Value = 3044.4868
$ws.Range("M32").Value = -2757.4868

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6078.7617
$ws.Range("J45").Value = 1640.25
$ws.Range("L45").Value = 1640.25
$ws.Range("N45").Value = -2394.25

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3081.0322
$ws.Range("I61").Value = 2560.52
$ws.Range("J61").Value = 5249.8335
$ws.Range("K61").Value = 2560.52
$ws.Range("L61").Value = 5249.8335
$ws.Range("M61").Value = -2348.52
$ws.Range("N61").Value = -5673.8335

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 8185.5264
$ws.Range("I74").Value = 1623.909
$ws.Range("J74").Value = 17207.75
$ws.Range("K74").Value = 1623.909
$ws.Range("L74").Value = 17207.75
$ws.Range("M74").Value = -749.9090000000001
$ws.Range("N74").Value = -18955.75

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 8185.5264
$ws.Range("I77").Value = 1623.909
$ws.Range("J77").Value = 17207.75
$ws.Range("K77").Value = 8119.545
$ws.Range("L77").Value = 86038.75
$ws.Range("M77").Value = -3751.545
$ws.Range("N77").Value = -94774.75

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1581.7097
$ws.Range("I97").Value = 1156.3103
$ws.Range("K97").Value = 1156.3103
$ws.Range("M97").Value = -660.3103000000001

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3549.5
$ws.Range("I102").Value = 3279.6
$ws.Range("K102").Value = 3279.6
$ws.Range("M102").Value = -1657.6

# ARM row 103
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = $null

# ARM row 112
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 30166.334
$ws.Range("J112").Value = 30166.334
$ws.Range("L112").Value = 30166.334
$ws.Range("N112").Value = -33120.334

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 478.72726
$ws.Range("I116").Value = 246.84616
$ws.Range("J116").Value = 813.6667
$ws.Range("K116").Value = 246.84616
$ws.Range("L116").Value = 813.6667
$ws.Range("M116").Value = 2047.15384
$ws.Range("N116").Value = -5401.6667

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3167.0952
$ws.Range("I132").Value = 2860.1516
$ws.Range("K132").Value = 8580.4548
$ws.Range("M132").Value = -6050.4548

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3081.0322
$ws.Range("I136").Value = 2560.52
$ws.Range("J136").Value = 5249.8335
$ws.Range("K136").Value = 7681.559999999999
$ws.Range("L136").Value = 15749.5005
$ws.Range("M136").Value = -5131.559999999999
$ws.Range("N136").Value = -20849.5005

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 478.72726
$ws.Range("I3").Value = 246.84616
$ws.Range("J3").Value = 813.6667
$ws.Range("K3").Value = 246.84616
$ws.Range("L3").Value = 813.6667
$ws.Range("M3").Value = -132.84616
$ws.Range("N3").Value = -1041.6667

# BSM row 74
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 79999
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = $null

# BSM row 77
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H77").Value = 79999
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = $null

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 207479.8
$ws.Range("I99").Value = 204701.6
$ws.Range("K99").Value = 204701.6
$ws.Range("M99").Value = -203203.6

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1359.1904
$ws.Range("I107").Value = 1148.8823
$ws.Range("K107").Value = 1148.8823
$ws.Range("M107").Value = 771.1177

# BSM row 138
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 71250
$ws.Range("J138").Value = 71250
$ws.Range("L138").Value = 71250
$ws.Range("N138").Value = -81530

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 464
$ws.Range("I22").Value = 246.83333
$ws.Range("J22").Value = 724.6
$ws.Range("K22").Value = 246.83333
$ws.Range("L22").Value = 724.6
$ws.Range("M22").Value = 103.16667
$ws.Range("N22").Value = -1424.6

# CRP row 26
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 499
$ws.Range("J26").Value = 499
$ws.Range("L26").Value = 499
$ws.Range("N26").Value = -1073

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 69787.266
$ws.Range("I31").Value = 127219.625
$ws.Range("J31").Value = 4150.2856
$ws.Range("K31").Value = 127219.625
$ws.Range("L31").Value = 4150.2856
$ws.Range("M31").Value = -126924.625
$ws.Range("N31").Value = -4740.2856

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 69787.266
$ws.Range("I34").Value = 127219.625
$ws.Range("J34").Value = 4150.2856
$ws.Range("K34").Value = 127219.625
$ws.Range("L34").Value = 4150.2856
$ws.Range("M34").Value = -127017.625
$ws.Range("N34").Value = -4554.2856

# CRP row 52
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 101281.8
$ws.Range("J52").Value = 104175
$ws.Range("L52").Value = 104175
$ws.Range("N52").Value = -104763

# CRP row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 18072.5

# CRP row 103
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 15047.833
$ws.Range("I103").Value = 10057.4
$ws.Range("J103").Value = 40000
$ws.Range("K103").Value = 10057.4
$ws.Range("L103").Value = 40000
$ws.Range("M103").Value = -8885.4
$ws.Range("N103").Value = -42344

# CRP row 124
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H124").Value = 46995
$ws.Range("J124").Value = 46995
$ws.Range("L124").Value = 46995
$ws.Range("N124").Value = -51905

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4509.533
$ws.Range("I132").Value = 4474.5713
$ws.Range("K132").Value = 13423.7139
$ws.Range("M132").Value = -10893.7139

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 13120.238
$ws.Range("I134").Value = 7501.5293
$ws.Range("J134").Value = 36999.75
$ws.Range("K134").Value = 22504.5879
$ws.Range("L134").Value = 110999.25
$ws.Range("M134").Value = -19969.5879
$ws.Range("N134").Value = -116069.25

# CRP row 135
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 65900
$ws.Range("J135").Value = 65900
$ws.Range("L135").Value = 65900
$ws.Range("N135").Value = -76040

# CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 15.6
$ws.Range("I2").Value = 19.095238
$ws.Range("J2").Value = 7.4444447
$ws.Range("K2").Value = 114.571428
$ws.Range("L2").Value = 44.6666682
$ws.Range("M2").Value = -1.571427999999997
$ws.Range("N2").Value = -270.6666682

# CUL row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 105.5

# CUL row 70
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 7771.143
$ws.Range("I70").Value = 6880
$ws.Range("K70").Value = 20640
$ws.Range("M70").Value = -20325

# CUL row 73
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 7771.143
$ws.Range("I73").Value = 6880
$ws.Range("K73").Value = 20640
$ws.Range("M73").Value = -19548

# CUL row 75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2770
$ws.Range("I75").Value = 1950
$ws.Range("K75").Value = 5850
$ws.Range("M75").Value = -4852

# CUL row 78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 2770
$ws.Range("I78").Value = 1950
$ws.Range("K78").Value = 17550
$ws.Range("M78").Value = -12558

# CUL row 100
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 9000
$ws.Range("J100").Value = 9000
$ws.Range("L100").Value = 27000
$ws.Range("N100").Value = -28622

# CUL row 102
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 7663.8335
$ws.Range("J102").Value = 7997.4116
$ws.Range("L102").Value = 23992.2348
$ws.Range("N102").Value = -28860.2348

# CUL row 104
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 9998.286
$ws.Range("J104").Value = 9998.286
$ws.Range("L104").Value = 29994.858
$ws.Range("N104").Value = -35236.858

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 636.3570999999999
$ws.Range("I107").Value = 526.2222
$ws.Range("J107").Value = 834.6
$ws.Range("K107").Value = 1578.6666
$ws.Range("L107").Value = 2503.8
$ws.Range("M107").Value = 341.3334
$ws.Range("N107").Value = -6343.8

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2016.6666
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2016.6666
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 18149.9994
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = -23049.9994

# CUL row 126
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 3015
$ws.Range("I126").Value = 3015
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9045
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4105
$ws.Range("N126").Value = $null

# GSM row 44
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = $null

# GSM row 52
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 50000
$ws.Range("I52").Value = 20000
$ws.Range("J52").Value = 60000
$ws.Range("K52").Value = 20000
$ws.Range("L52").Value = 60000
$ws.Range("M52").Value = -19741
$ws.Range("N52").Value = -60518

# GSM row 53
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").Value = $null

# GSM row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 7000
$ws.Range("I57").Value = 7000
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 7000
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -6180
$ws.Range("N57").Value = $null

# GSM row 58
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 60000
$ws.Range("J58").Value = 60000
$ws.Range("L58").Value = 60000
$ws.Range("N58").Value = -60554

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3172.5557
$ws.Range("I80").Value = 3108
$ws.Range("J80").Value = 3301.6667
$ws.Range("K80").Value = 3108
$ws.Range("L80").Value = 3301.6667
$ws.Range("M80").Value = -2110
$ws.Range("N80").Value = -5297.6667

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3172.5557
$ws.Range("I83").Value = 3108
$ws.Range("J83").Value = 3301.6667
$ws.Range("K83").Value = 15540
$ws.Range("L83").Value = 16508.3335
$ws.Range("M83").Value = -10548
$ws.Range("N83").Value = -26492.3335

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2910.75
$ws.Range("I122").Value = 2307.5217
$ws.Range("J122").Value = 4452.3335
$ws.Range("K122").Value = 6922.5651
$ws.Range("L122").Value = 13357.0005
$ws.Range("M122").Value = -4472.5651
$ws.Range("N122").Value = -18257.0005

# LTW row 9
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 527
$ws.Range("I9").Value = 399
$ws.Range("K9").Value = 399
$ws.Range("M9").Value = -175

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2569.3928
$ws.Range("I16").Value = 1920.8846
$ws.Range("J16").Value = 11000
$ws.Range("K16").Value = 1920.8846
$ws.Range("L16").Value = 11000
$ws.Range("M16").Value = -1750.8846
$ws.Range("N16").Value = -11340

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1783
$ws.Range("I22").Value = 1036.75
$ws.Range("K22").Value = 1036.75
$ws.Range("M22").Value = -741.75

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1783
$ws.Range("I27").Value = 1036.75
$ws.Range("K27").Value = 1036.75
$ws.Range("M27").Value = -929.75

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3832.7083
$ws.Range("I40").Value = 3498.6667
$ws.Range("J40").Value = 4389.4443
$ws.Range("K40").Value = 3498.6667
$ws.Range("L40").Value = 4389.4443
$ws.Range("M40").Value = -3362.6667
$ws.Range("N40").Value = -4661.4443

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2999.875
$ws.Range("I68").Value = 2999.8572
$ws.Range("K68").Value = 2999.8572
$ws.Range("M68").Value = -2250.8572

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2999.875
$ws.Range("I71").Value = 2999.8572
$ws.Range("K71").Value = 14999.286
$ws.Range("M71").Value = -11255.286

# LTW row 81
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 10082.083
$ws.Range("I82").Value = 14091.5
$ws.Range("K82").Value = 14091.5
$ws.Range("M82").Value = -13730.5

# LTW row 84
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 10082.083
$ws.Range("I85").Value = 14091.5
$ws.Range("K85").Value = 14091.5
$ws.Range("M85").Value = -12843.5

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2024.5
$ws.Range("I100").Value = 1549
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 1549
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -1008
$ws.Range("N100").Value = -3582

# LTW row 138
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H138").Value = 67058
$ws.Range("I138").Value = 50390
$ws.Range("J138").Value = 71225
$ws.Range("K138").Value = 50390
$ws.Range("L138").Value = 71225
$ws.Range("M138").Value = -45250
$ws.Range("N138").Value = -81505

# LTW row 139
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 96602
$ws.Range("J139").Value = 96602
$ws.Range("L139").Value = 96602
$ws.Range("N139").Value = -106882

# WVR row 14
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 3665.6667
$ws.Range("J14").Value = 3665.6667
$ws.Range("L14").Value = 3665.6667
$ws.Range("N14").Value = -4001.6667

# WVR row 19
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").Value = $null

# WVR row 64
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null

# WVR row 67
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 43480.832
$ws.Range("I136").Value = 67929.8
$ws.Range("K136").Value = 203789.4
$ws.Range("M136").Value = -201239.4
